# Applies the "Add Bomb, score, pause" edit to the scoring worksheet:
#  - adds a new shared string "쓸꺼임!" used as a note on row 16 (D/E columns)
#  - adds several new "D" column mirror-values for rows that previously
#    only had a C column score (16, 22, 41-45)
#  - removes the stray "예정" note that used to sit in E29
#  - bumps the score review row 50's D value from 3 to 5
#  - moves the saved scroll/selection position of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheetData edits -------------------------------------------------

# Row 16: add a "D" score mirror plus a little sticky-note comment in E.
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = "쓸꺼임!"

# Row 22: add the "D" score mirror (no note).
$ws.Range("D22").Value = 1

# Row 29: the old "예정" (pending) note is no longer needed.
$ws.Range("E29").ClearContents()

# Rows 41-45: add "D" score mirrors matching the existing "C" values.
$ws.Range("D41").Value = 3
$ws.Range("D42").Value = 5
$ws.Range("D43").Value = 8
$ws.Range("D44").Value = 4
$ws.Range("D45").Value = 5

# Row 50: bump the awarded score from 3 up to 5.
$ws.Range("D50").Value = 5

# --- view state --------------------------------------------------------
# Scroll the visible window down a bit and move the active selection,
# same as what a reviewer scrolling further down the sheet would do.
$ws.Range("M30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
